$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter => new value (same value applies to rows 2 and 3)
$updates = [ordered]@{
    "G"  = 0.15625
    "H"  = 0.15625
    "I"  = 0.2449404761904762
    "J"  = 0.1794226555652936
    "K"  = 5.97
    "L"  = 0.1776785714285714
    "M"  = 2
    "N"  = 0.066006600660066
    "O"  = 0.3350083752093803
    "P"  = 2
    "Q"  = 0.066006600660066
    "R"  = 0.3350083752093803
    "U"  = 0.554
    "V"  = 0.01828382838283828
    "W"  = 0.312565445026178
    "X"  = 0.08860093357074675
    "Y"  = 0.2239645114554313
    "Z"  = 1.748815905897049
    "AA" = 0.3137771939308731
    "AB" = 0.08819342929913308
    "AC" = 0.22558376463174
    "AD" = 0.463
    "AF" = 0.463
    "AG" = -0.09100000000000003
    "AH" = 0.01505054773591652
    "AI" = 0.01924115862527532
    "AJ" = -0.003012347313714457
    "AK" = -0.003870857969288358
    "AL" = 0.111
    "AM" = 0.111
    "AN" = 0.05518474374255065
    "AO" = 74.14414414414415
    "AP" = -0.01084624553039333
    "AQ" = 74.14414414414415
}

foreach ($col in $updates.Keys) {
    $value = $updates[$col]
    $ws.Range("$col" + "2").Value = $value
    $ws.Range("$col" + "3").Value = $value
}
